# Fruta / hortaliza, semanal
# Insert two new weekly observation rows (new rows 215 and 216) into the
# Frambuesa / Lo Valledor dataset, pushing the existing rows 215-237 down
# to 217-239.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 215, shifting the
# remaining data (old rows 215-237) down to rows 217-239.
$ws.Rows.Item(215).Insert()
$ws.Rows.Item(216).Insert()

# --- New row 215 ---
$ws.Range("A215").Value = 6
$ws.Range("B215").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C215").Value = "Metropolitana"
$ws.Range("D215").Value = 44918
$ws.Range("D215").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E215").Value = 13
$ws.Range("F215").Value = "Fruta"
$ws.Range("G215").Value = 100101
$ws.Range("H215").Value = "Berries"
$ws.Range("I215").Value = 100101004
$ws.Range("J215").Value = "Frambuesa"
$ws.Range("K215").Value = "Sin especificar"
$ws.Range("L215").Value = "Especial"
$ws.Range("M215").Value = 300
$ws.Range("N215").Value = 7000
$ws.Range("O215").Value = 7000
$ws.Range("P215").Value = 7000
$ws.Range("Q215").Value = "`$/bandeja 2 kilos"
$ws.Range("R215").Value = "Provincia de Curicó"
$ws.Range("S215").Value = 3500
$ws.Range("T215").Value = 2

# --- New row 216 ---
$ws.Range("A216").Value = 6
$ws.Range("B216").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C216").Value = "Metropolitana"
$ws.Range("D216").Value = 44918
$ws.Range("D216").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E216").Value = 13
$ws.Range("F216").Value = "Fruta"
$ws.Range("G216").Value = 100101
$ws.Range("H216").Value = "Berries"
$ws.Range("I216").Value = 100101004
$ws.Range("J216").Value = "Frambuesa"
$ws.Range("K216").Value = "Sin especificar"
$ws.Range("L216").Value = "Especial"
$ws.Range("M216").Value = 250
$ws.Range("N216").Value = 7000
$ws.Range("O216").Value = 7000
$ws.Range("P216").Value = 7000
$ws.Range("Q216").Value = "`$/bandeja 2 kilos"
$ws.Range("R216").Value = "Región del Maule"
$ws.Range("S216").Value = 3500
$ws.Range("T216").Value = 2

$ws.Range("A1").Select()

Write-Host "Final used range:" $ws.UsedRange.Address()
